$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.141788666666667
$ws.Range("H2").Value = 3.425366
$ws.Range("I2").Value = 0.2235063054668694
$ws.Range("J2").Value = 0.2235063054668694
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5203476666666668
$ws.Range("N2").Value = 1.561043
$ws.Range("O2").Value = 0.004105934376266647
$ws.Range("P2").Value = 0.004105934376266647
$ws.Range("Q2").Value = 0.5941270685264445
$ws.Range("R2").Value = 5.347143616738
$ws.Range("S2").Value = 0.0009177022229287732
$ws.Range("T2").Value = 0.0009177022229287733

$ws.Range("G3").Value = 1.141788666666667
$ws.Range("H3").Value = 3.425366
$ws.Range("I3").Value = 0.2235063054668694
$ws.Range("J3").Value = 0.2235063054668694
$ws.Range("O3").Value = 0.8361295370252257
$ws.Range("P3").Value = 0.8361295370252259
$ws.Range("Q3").Value = 120.9876109108349
$ws.Range("R3").Value = 1088.888498197514
$ws.Range("S3").Value = 0.1868802237122322
$ws.Range("T3").Value = 0.1868802237122322

$ws.Range("G4").Value = 1.141788666666667
$ws.Range("H4").Value = 3.425366
$ws.Range("I4").Value = 0.2235063054668694
$ws.Range("J4").Value = 0.2235063054668694
$ws.Range("M4").Value = 20.24706
$ws.Range("N4").Value = 60.74118
$ws.Range("O4").Value = 0.1597645285985076
$ws.Range("P4").Value = 0.1597645285985076
$ws.Range("Q4").Value = 23.11786364132
$ws.Range("R4").Value = 208.06077277188
$ws.Range("S4").Value = 0.03570837953170844
$ws.Range("T4").Value = 0.03570837953170844

$ws.Range("I5").Value = 0.5880650598117267
$ws.Range("J5").Value = 0.5880650598117266
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5203476666666668
$ws.Range("N5").Value = 1.561043
$ws.Range("O5").Value = 0.004105934376266647
$ws.Range("P5").Value = 0.004105934376266647
$ws.Range("Q5").Value = 1.563201402121334
$ws.Range("R5").Value = 14.068812619092
$ws.Range("S5").Value = 0.00241455654456227
$ws.Range("T5").Value = 0.002414556544562271

$ws.Range("I6").Value = 0.5880650598117267
$ws.Range("J6").Value = 0.5880650598117266
$ws.Range("O6").Value = 0.8361295370252257
$ws.Range("P6").Value = 0.8361295370252259
$ws.Range("S6").Value = 0.4916985662010908
$ws.Range("T6").Value = 0.4916985662010908

$ws.Range("I7").Value = 0.5880650598117267
$ws.Range("J7").Value = 0.5880650598117266
$ws.Range("M7").Value = 20.24706
$ws.Range("N7").Value = 60.74118
$ws.Range("O7").Value = 0.1597645285985076
$ws.Range("P7").Value = 0.1597645285985076
$ws.Range("Q7").Value = 60.82516480488001
$ws.Range("R7").Value = 547.4264832439201
$ws.Range("S7").Value = 0.0939519370660737
$ws.Range("T7").Value = 0.09395193706607369

$ws.Range("G8").Value = 0.9625933333333334
$ws.Range("H8").Value = 2.88778
$ws.Range("I8").Value = 0.1884286347214039
$ws.Range("J8").Value = 0.1884286347214039
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5203476666666668
$ws.Range("N8").Value = 1.561043
$ws.Range("O8").Value = 0.004105934376266647
$ws.Range("P8").Value = 0.004105934376266647
$ws.Range("Q8").Value = 0.500883194948889
$ws.Range("R8").Value = 4.507948754540001
$ws.Range("S8").Value = 0.0007736756087756033
$ws.Range("T8").Value = 0.0007736756087756033

$ws.Range("G9").Value = 0.9625933333333334
$ws.Range("H9").Value = 2.88778
$ws.Range("I9").Value = 0.1884286347214039
$ws.Range("J9").Value = 0.1884286347214039
$ws.Range("O9").Value = 0.8361295370252257
$ws.Range("P9").Value = 0.8361295370252259
$ws.Range("Q9").Value = 101.9994952469578
$ws.Range("R9").Value = 917.9954572226201
$ws.Range("S9").Value = 0.1575507471119028
$ws.Range("T9").Value = 0.1575507471119028

$ws.Range("G10").Value = 0.9625933333333334
$ws.Range("H10").Value = 2.88778
$ws.Range("I10").Value = 0.1884286347214039
$ws.Range("J10").Value = 0.1884286347214039
$ws.Range("M10").Value = 20.24706
$ws.Range("N10").Value = 60.74118
$ws.Range("O10").Value = 0.1597645285985076
$ws.Range("P10").Value = 0.1597645285985076
$ws.Range("Q10").Value = 19.4896849756
$ws.Range("R10").Value = 175.4071647804
$ws.Range("S10").Value = 0.03010421200072547
$ws.Range("T10").Value = 0.03010421200072547
